$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 528023.4
$ws.Range("J17").Value = 528023.4
$ws.Range("L17").Value = 1584070.2
$ws.Range("N17").Value = -1584406.2
$ws.Range("H111").Value = 544.125
$ws.Range("I111").Value = 394.5
$ws.Range("K111").Value = 1183.5
$ws.Range("M111").Value = 1883.5
$ws.Range("H115").Value = 28890702
$ws.Range("I115").Value = 32501814
$ws.Range("J115").Value = 1800
$ws.Range("K115").Value = 97505442
$ws.Range("L115").Value = 5400
$ws.Range("M115").Value = -97503875
$ws.Range("N115").Value = -8534
$ws.Range("H138").Value = 2656.9714
$ws.Range("J138").Value = 3203.0476
$ws.Range("L138").Value = 9609.1428
$ws.Range("N138").Value = -19889.1428

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2544.5862
$ws.Range("I2").Value = 2334.9614
$ws.Range("J2").Value = 4361.3335
$ws.Range("K2").Value = 2334.9614
$ws.Range("L2").Value = 4361.3335
$ws.Range("M2").Value = -2221.9614
$ws.Range("N2").Value = -4587.3335
$ws.Range("H32").Value = 4074.0547
$ws.Range("I32").Value = 3505.8142
$ws.Range("K32").Value = 3505.8142
$ws.Range("M32").Value = -3218.8142
$ws.Range("H45").Value = 5031.0527
$ws.Range("J45").Value = 3669.6875
$ws.Range("L45").Value = 3669.6875
$ws.Range("N45").Value = -4423.6875
$ws.Range("H61").Value = 1642.6666
$ws.Range("I61").Value = 1559.2222
$ws.Range("K61").Value = 1559.2222
$ws.Range("M61").Value = -1347.2222
$ws.Range("H74").Value = 2454.158
$ws.Range("I74").Value = 2098
$ws.Range("J74").Value = 2943.875
$ws.Range("K74").Value = 2098
$ws.Range("L74").Value = 2943.875
$ws.Range("M74").Value = -1224
$ws.Range("N74").Value = -4691.875
$ws.Range("H77").Value = 2454.158
$ws.Range("I77").Value = 2098
$ws.Range("J77").Value = 2943.875
$ws.Range("K77").Value = 10490
$ws.Range("L77").Value = 14719.375
$ws.Range("M77").Value = -6122
$ws.Range("N77").Value = -23455.375
$ws.Range("H116").Value = 2544.5862
$ws.Range("I116").Value = 2334.9614
$ws.Range("J116").Value = 4361.3335
$ws.Range("K116").Value = 2334.9614
$ws.Range("L116").Value = 4361.3335
$ws.Range("M116").Value = -40.96140000000014
$ws.Range("N116").Value = -8949.333500000001
$ws.Range("H122").Value = 2521.2258
$ws.Range("J122").Value = 3400.4167
$ws.Range("L122").Value = 10201.2501
$ws.Range("N122").Value = -15101.2501
$ws.Range("H132").Value = 10649.818
$ws.Range("I132").Value = 10649.818
$ws.Range("K132").Value = 31949.454
$ws.Range("M132").Value = -29419.454
$ws.Range("H136").Value = 1642.6666
$ws.Range("I136").Value = 1559.2222
$ws.Range("K136").Value = 4677.6666
$ws.Range("M136").Value = -2127.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2544.5862
$ws.Range("I3").Value = 2334.9614
$ws.Range("J3").Value = 4361.3335
$ws.Range("K3").Value = 2334.9614
$ws.Range("L3").Value = 4361.3335
$ws.Range("M3").Value = -2220.9614
$ws.Range("N3").Value = -4589.3335
$ws.Range("H132").Value = 134783.5
$ws.Range("J132").Value = 134783.5
$ws.Range("L132").Value = 134783.5
$ws.Range("N132").Value = -144903.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 25634.512
$ws.Range("I31").Value = 31483
$ws.Range("J31").Value = 3540.2222
$ws.Range("K31").Value = 31483
$ws.Range("L31").Value = 3540.2222
$ws.Range("M31").Value = -31188
$ws.Range("N31").Value = -4130.2222
$ws.Range("H34").Value = 25634.512
$ws.Range("I34").Value = 31483
$ws.Range("J34").Value = 3540.2222
$ws.Range("K34").Value = 31483
$ws.Range("L34").Value = 3540.2222
$ws.Range("M34").Value = -31281
$ws.Range("N34").Value = -3944.2222
$ws.Range("H58").Value = 1849
$ws.Range("I58").Value = 1868
$ws.Range("J58").Value = 1798.3334
$ws.Range("K58").Value = 1868
$ws.Range("L58").Value = 1798.3334
$ws.Range("M58").Value = -1665
$ws.Range("N58").Value = -2204.3334
$ws.Range("H132").Value = 3693.5557
$ws.Range("I132").Value = 3964.3914
$ws.Range("K132").Value = 11893.1742
$ws.Range("M132").Value = -9363.174199999999
$ws.Range("H136").Value = 1849
$ws.Range("I136").Value = 1868
$ws.Range("J136").Value = 1798.3334
$ws.Range("K136").Value = 5604
$ws.Range("L136").Value = 5395.0002
$ws.Range("M136").Value = -3054
$ws.Range("N136").Value = -10495.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 610.0625
$ws.Range("I122").Value = 636.2
$ws.Range("J122").Value = 566.5
$ws.Range("K122").Value = 5725.8
$ws.Range("L122").Value = 5098.5
$ws.Range("M122").Value = -3275.8
$ws.Range("N122").Value = -9998.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2516.8262
$ws.Range("I102").Value = 2682.2222
$ws.Range("K102").Value = 2682.2222
$ws.Range("M102").Value = -1060.2222
$ws.Range("H140").Value = 117358.555
$ws.Range("J140").Value = 117358.555
$ws.Range("L140").Value = 117358.555
$ws.Range("N140").Value = -127718.555

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 19175.133
$ws.Range("I7").Value = 31150.555
$ws.Range("K7").Value = 31150.555
$ws.Range("M7").Value = -31038.555
$ws.Range("H22").Value = 834.2727
$ws.Range("I22").Value = 691.8570999999999
$ws.Range("J22").Value = 900.73334
$ws.Range("K22").Value = 691.8570999999999
$ws.Range("L22").Value = 900.73334
$ws.Range("M22").Value = -396.8570999999999
$ws.Range("N22").Value = -1490.73334
$ws.Range("H27").Value = 834.2727
$ws.Range("I27").Value = 691.8570999999999
$ws.Range("J27").Value = 900.73334
$ws.Range("K27").Value = 691.8570999999999
$ws.Range("L27").Value = 900.73334
$ws.Range("M27").Value = -584.8570999999999
$ws.Range("N27").Value = -1114.73334
$ws.Range("H40").Value = 67633
$ws.Range("I40").Value = 22383.166
$ws.Range("K40").Value = 22383.166
$ws.Range("M40").Value = -22247.166
$ws.Range("H61").Value = 80596.234
$ws.Range("I61").Value = 78093.46000000001
$ws.Range("K61").Value = 78093.46000000001
$ws.Range("M61").Value = -77891.46000000001
$ws.Range("H68").Value = 2293.9167
$ws.Range("I68").Value = 2472.7
$ws.Range("J68").Value = 1400
$ws.Range("K68").Value = 2472.7
$ws.Range("L68").Value = 1400
$ws.Range("M68").Value = -1723.7
$ws.Range("N68").Value = -2898
$ws.Range("H71").Value = 2293.9167
$ws.Range("I71").Value = 2472.7
$ws.Range("J71").Value = 1400
$ws.Range("K71").Value = 12363.5
$ws.Range("L71").Value = 7000
$ws.Range("M71").Value = -8619.5
$ws.Range("N71").Value = -14488
$ws.Range("H87").Value = 22594.5
$ws.Range("H88").Value = 34404
$ws.Range("H90").Value = 22594.5
$ws.Range("H91").Value = 34404
$ws.Range("H113").Value = 80596.234
$ws.Range("I113").Value = 78093.46000000001
$ws.Range("K113").Value = 78093.46000000001
$ws.Range("M113").Value = -75923.46000000001
$ws.Range("H126").Value = 19175.133
$ws.Range("I126").Value = 31150.555
$ws.Range("K126").Value = 93451.66500000001
$ws.Range("M126").Value = -90981.66500000001
$ws.Range("H132").Value = 3824.225
$ws.Range("I132").Value = 3203.6897
$ws.Range("J132").Value = 5460.1816
$ws.Range("K132").Value = 9611.069100000001
$ws.Range("L132").Value = 16380.5448
$ws.Range("M132").Value = -7081.069100000001
$ws.Range("N132").Value = -21440.5448

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2596.6667
$ws.Range("I81").Value = 2596.6667
$ws.Range("K81").Value = 5193.3334
$ws.Range("M81").Value = -4132.3334
$ws.Range("H84").Value = 2596.6667
$ws.Range("I84").Value = 2596.6667
$ws.Range("K84").Value = 25966.667
$ws.Range("M84").Value = -20662.667
